# Updates the cryptos price/volume(1h) table with the latest scraped values.
# Price-column (D) values are assigned with a leading "'" (quote-prefix) so
# Excel keeps them as literal text instead of re-parsing numeric-looking
# strings (e.g. "19.90", "0.0470") into numbers and dropping trailing/
# leading zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.025.04'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '''3.777.21'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''429.36'
$ws.Range('E5').Value = '  +5.20%  '
$ws.Range('D6').Value = '''138.49'
$ws.Range('E6').Value = '  +4.35%  '
$ws.Range('D7').Value = '''0.621'
$ws.Range('E7').Value = '  +2.08%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').Value = '''0.735'
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('E10').Value = '  -9.62%  '
$ws.Range('D11').Value = '''0.0000309'
$ws.Range('E11').Value = '  -13.82%  '
$ws.Range('D12').Value = '''42.62'
$ws.Range('E12').Value = '  +3.82%  '
$ws.Range('D13').Value = '''10.42'
$ws.Range('E13').Value = '  +5.41%  '
$ws.Range('D14').Value = '''4.375.10'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '''3.782.86'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '''19.90'
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('E19').Value = '  +5.57%  '
$ws.Range('D20').Value = '''66.114.61'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').Value = '''404.16'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('D22').Value = '''14.88'
$ws.Range('E22').Value = '  +3.44%  '
$ws.Range('E23').Value = '  +7.38%  '
$ws.Range('D24').Value = '''84.53'
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').Value = '''10.09'
$ws.Range('E25').Value = '  +36.18%  '
$ws.Range('D26').Value = '''36.52'
$ws.Range('E26').Value = '  +1.45%  '
$ws.Range('D27').Value = '''3.30'
$ws.Range('E27').Value = '  +5.67%  '
$ws.Range('E28').Value = '  -3.84%  '
$ws.Range('D29').Value = '''9.83'
$ws.Range('E29').Value = '  +5.01%  '
$ws.Range('E30').Value = '  +13.51%  '
$ws.Range('D31').Value = '''13.73'
$ws.Range('E31').Value = '  +11.08%  '
$ws.Range('D32').Value = '''705.33'
$ws.Range('E32').Value = '  -4.68%  '
$ws.Range('E33').Value = '  +3.36%  '
$ws.Range('D34').Value = '''41.16'
$ws.Range('E34').Value = '  +5.35%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '''5.77'
$ws.Range('E36').Value = '  +35.60%  '
$ws.Range('E37').Value = '  -3.38%  '
$ws.Range('D38').Value = '''56.12'
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0470'
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '''2.84'
$ws.Range('E40').Value = '  +41.76%  '
$ws.Range('D41').Value = '''3.02'
$ws.Range('E41').Value = '  +6.22%  '
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('D44').Value = '''0.0₃0667'
$ws.Range('E44').Value = '  -11.54%  '
$ws.Range('D45').Value = '''0.330'
$ws.Range('E45').Value = '  +11.52%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.19'
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').Value = '''3.35'
$ws.Range('E47').Value = '  +2.41%  '
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').Value = '''138.84'
$ws.Range('E50').Value = '  -4.02%  '
$ws.Range('E51').Value = '  -0.62%  '
